$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New test-case rows (21-25) -------------------------------------------
# Column B (test-case names) is populated first, then column A (category
# labels), matching the order new shared strings were introduced by the
# original author's edit.
$ws.Range("B21").Value = "011_Storage_Type_Add_Holds_Multiple_Container"
$ws.Range("B22").Value = "08_Storage_Type_Verify_Default_Temp"
$ws.Range("B23").Value = "09_Storage_Type_Error_Temperature"
$ws.Range("B25").Value = "04_Storage_Type_Edit_Error"
$ws.Range("B24").Value = "02_Storage_Type_Add_Error"

$ws.Range("A21").Value = "Storage_Type_Add_success"
$ws.Range("A22").Value = "Storage_Type_Add_success"
$ws.Range("A23").Value = "Storage_Type_Edit_error"
$ws.Range("A24").Value = "Storage_Type_Add_error"
$ws.Range("A25").Value = "Storage_type_edit_error"

# Column A uses the bold "category" style already used elsewhere in the sheet
$ws.Range("A21:A25").Font.Bold = $true

# --- Header row grows to a two-line height ---------------------------------
$ws.Rows(1).RowHeight = 30

# --- Selection moves onto the newly added rows ------------------------------
$ws.Range("A21:A25").Select() | Out-Null
